$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing June day-13 total_venda value
$ws.Cells.Item(14, 2).Value = 22435.29

# Insert 2 new rows for June days 14 and 15 (this shifts all subsequent
# rows down by 2, which automatically relocates the May/April/March data
# into their correct final rows, including the 2 "new" March rows that
# appear at the tail of the sheet in the diff - those are just the
# pre-existing March 29/30 rows landing in their shifted positions).
$ws.Rows("15:16").Insert()

# Fill in the two newly inserted June rows
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 27509.4
$ws.Cells.Item(15, 3).Value = 6
$ws.Cells.Item(15, 4).Value = 2025
$ws.Cells.Item(15, 5).Value = "06/2025"

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 4621.42
$ws.Cells.Item(16, 3).Value = 6
$ws.Cells.Item(16, 4).Value = 2025
$ws.Cells.Item(16, 5).Value = "06/2025"
